$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Group Size 3")

# Row 2: Alice=0.3, Bob=0.7, Carol stays 0
$ws.Range("D2").Value = 0.3
$ws.Range("E2").Value = 0.7

# Row 3: Alice=1
$ws.Range("D3").Value = 1

# Row 4: Carol=1
$ws.Range("F4").Value = 1

# Row 6: Carol=1
$ws.Range("F6").Value = 1

# Row 7: Alice=1
$ws.Range("D7").Value = 1

# Row 8: Alice=1
$ws.Range("D8").Value = 1

# Row 10: Alice=1
$ws.Range("D10").Value = 1

# Row 11: Alice=1
$ws.Range("D11").Value = 1

# Row 12: Alice=1
$ws.Range("D12").Value = 1

# Row 13: Alice=1
$ws.Range("D13").Value = 1

# Row 14: Alice=1
$ws.Range("D14").Value = 1

# Row 15: Alice=0.4, Bob=0.3, Carol=0.3
$ws.Range("D15").Value = 0.4
$ws.Range("E15").Value = 0.3
$ws.Range("F15").Value = 0.3

# Row 16: Alice=0.4, Bob=0.3, Carol=0.3
$ws.Range("D16").Value = 0.4
$ws.Range("E16").Value = 0.3
$ws.Range("F16").Value = 0.3

# Row 21: Alice=0.4, Bob=0.3, Carol=0.3
$ws.Range("D21").Value = 0.4
$ws.Range("E21").Value = 0.3
$ws.Range("F21").Value = 0.3

# Row 22: Alice=0.4, Bob=0.3, Carol=0.3
$ws.Range("D22").Value = 0.4
$ws.Range("E22").Value = 0.3
$ws.Range("F22").Value = 0.3

# Update selection to K15 to match the saved cursor position
$ws.Range("K15").Select()
